# Commit: "Removed the pseudo SQL function and added a join to Geography to the system prompt"
#
# 1) Rename the existing (only) worksheet from "Sheet1" to "GC7_budget".
# 2) Update its "cost_group"/"cost_input" pseudo-SQL-function header labels
#    (columns E/F, row 1) to human readable "Cost Category"/"Cost Input".
# 3) Add a new "Geography" worksheet right after GC7_budget containing a
#    Geography-name -> Region/Department join/lookup table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Fix up the pseudo "cost_group" / "cost_input" column headers on the budget
# sheet so they read as normal text instead of function-looking identifiers.
$ws1.Range("E1").Value = "Cost Category"
$ws1.Range("F1").Value = "Cost Input"

# Tidy up the view (no more pinned top-left / stray selection highlighting).
$ws1.Range("A1").Select()

# Rename Sheet1 -> GC7_budget
$ws1.Name = "GC7_budget"

# Add the new Geography sheet immediately after GC7_budget.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Geography"

# ---- Build the Geography sheet data table (country -> region -> dept reference) ----
$geoData = New-Object 'object[,]' 136,3
$geoData[0,0] = '[Geography Name]'
$geoData[0,1] = '[NewRegioShort]'
$geoData[0,2] = '[NewDept]'
$geoData[1,0] = 'Afghanistan'
$geoData[1,1] = 'Asia'
$geoData[1,2] = 'ALE'
$geoData[2,0] = 'Albania'
$geoData[2,1] = 'EECA'
$geoData[2,2] = 'AME'
$geoData[3,0] = 'Algeria'
$geoData[3,1] = 'HIA1'
$geoData[3,2] = 'HIA1'
$geoData[4,0] = 'Angola'
$geoData[4,1] = 'HIA2'
$geoData[4,2] = 'HIA2'
$geoData[5,0] = 'Armenia'
$geoData[5,1] = 'LAC'
$geoData[6,0] = 'Azerbaijan'
$geoData[6,1] = 'MENASEA'
$geoData[7,0] = 'Bangladesh'
$geoData[7,1] = 'WCA'
$geoData[8,0] = 'Belarus'
$geoData[9,0] = 'Belize'
$geoData[10,0] = 'Benin'
$geoData[11,0] = 'Bhutan'
$geoData[12,0] = 'Bolivia (Plurinational State)'
$geoData[13,0] = 'Botswana'
$geoData[14,0] = 'Burkina Faso'
$geoData[15,0] = 'Burundi'
$geoData[16,0] = 'Cabo Verde'
$geoData[17,0] = 'Cambodia'
$geoData[18,0] = 'Cameroon'
$geoData[19,0] = 'Central African Republic'
$geoData[20,0] = 'Chad'
$geoData[21,0] = 'Colombia'
$geoData[22,0] = 'Comoros'
$geoData[23,0] = 'Congo'
$geoData[24,0] = 'Congo (Democratic Republic)'
$geoData[25,0] = 'Costa Rica'
$geoData[26,0] = 'Côte d''Ivoire'
$geoData[27,0] = 'Cuba'
$geoData[28,0] = 'Djibouti'
$geoData[29,0] = 'Dominican Republic'
$geoData[30,0] = 'Ecuador'
$geoData[31,0] = 'Egypt'
$geoData[32,0] = 'El Salvador'
$geoData[33,0] = 'Equatorial Guinea'
$geoData[34,0] = 'Eritrea'
$geoData[35,0] = 'Eswatini'
$geoData[36,0] = 'Ethiopia'
$geoData[37,0] = 'Gabon'
$geoData[38,0] = 'Gambia'
$geoData[39,0] = 'Georgia'
$geoData[40,0] = 'Ghana'
$geoData[41,0] = 'Guatemala'
$geoData[42,0] = 'Guinea'
$geoData[43,0] = 'Guinea-Bissau'
$geoData[44,0] = 'Guyana'
$geoData[45,0] = 'Haiti'
$geoData[46,0] = 'Honduras'
$geoData[47,0] = 'India'
$geoData[48,0] = 'Indonesia'
$geoData[49,0] = 'Iran (Islamic Republic)'
$geoData[50,0] = 'Jamaica'
$geoData[51,0] = 'Kazakhstan'
$geoData[52,0] = 'Kenya'
$geoData[53,0] = 'Korea (Democratic Peoples Republic)'
$geoData[54,0] = 'Kosovo'
$geoData[55,0] = 'Kyrgyzstan'
$geoData[56,0] = 'Lao (Peoples Democratic Republic)'
$geoData[57,0] = 'Lesotho'
$geoData[58,0] = 'Liberia'
$geoData[59,0] = 'Madagascar'
$geoData[60,0] = 'Malawi'
$geoData[61,0] = 'Malaysia'
$geoData[62,0] = 'Mali'
$geoData[63,0] = 'Mauritania'
$geoData[64,0] = 'Mauritius'
$geoData[65,0] = 'Moldova'
$geoData[66,0] = 'Mongolia'
$geoData[67,0] = 'Montenegro'
$geoData[68,0] = 'Morocco'
$geoData[69,0] = 'Mozambique'
$geoData[70,0] = 'Multicountry Africa ECSA-HC'
$geoData[71,0] = 'Multicountry Americas EMMIE'
$geoData[72,0] = 'Multicountry Americas ORAS-CONHU'
$geoData[73,0] = 'Multicountry Caribbean CARICOM-PANCAP'
$geoData[74,0] = 'Multicountry Caribbean MCC'
$geoData[75,0] = 'Multicountry East Asia and Pacific RAI'
$geoData[76,0] = 'Multicountry Eastern Africa IGAD'
$geoData[77,0] = 'Multicountry EECA PAS'
$geoData[78,0] = 'Multicountry HIV EECA APH'
$geoData[79,0] = 'Multicountry HIV Latin America ALEP'
$geoData[80,0] = 'Multicountry HIV MENA IHAA'
$geoData[81,0] = 'Multicountry HIV SEA AFAO'
$geoData[82,0] = 'Multicountry MENA Key Populations'
$geoData[83,0] = 'Multicountry Middle East MER'
$geoData[84,0] = 'Multicountry Southern Africa E8'
$geoData[85,0] = 'Multicountry Southern Africa MOSASWA'
$geoData[86,0] = 'Multicountry Southern Africa TIMS'
$geoData[87,0] = 'Multicountry Southern Africa WHC'
$geoData[88,0] = 'Multicountry TB Asia TEAM'
$geoData[89,0] = 'Multicountry TB Asia UNDP'
$geoData[90,0] = 'Multicountry TB Asia UNOPS'
$geoData[91,0] = 'Multicountry TB LAC PIH'
$geoData[92,0] = 'Multicountry TB WC Africa NTP/SRL'
$geoData[93,0] = 'Multicountry West Africa ALCO'
$geoData[94,0] = 'Multicountry Western Pacific'
$geoData[95,0] = 'Myanmar'
$geoData[96,0] = 'Namibia'
$geoData[97,0] = 'Nepal'
$geoData[98,0] = 'Nicaragua'
$geoData[99,0] = 'Niger'
$geoData[100,0] = 'Nigeria'
$geoData[101,0] = 'Pakistan'
$geoData[102,0] = 'Panama'
$geoData[103,0] = 'Papua New Guinea'
$geoData[104,0] = 'Paraguay'
$geoData[105,0] = 'Peru'
$geoData[106,0] = 'Philippines'
$geoData[107,0] = 'Romania'
$geoData[108,0] = 'Russian Federation'
$geoData[109,0] = 'Rwanda'
$geoData[110,0] = 'Sao Tome and Principe'
$geoData[111,0] = 'Senegal'
$geoData[112,0] = 'Serbia'
$geoData[113,0] = 'Sierra Leone'
$geoData[114,0] = 'Solomon Islands'
$geoData[115,0] = 'Somalia'
$geoData[116,0] = 'South Africa'
$geoData[117,0] = 'South Sudan'
$geoData[118,0] = 'Sri Lanka'
$geoData[119,0] = 'Sudan'
$geoData[120,0] = 'Suriname'
$geoData[121,0] = 'Tajikistan'
$geoData[122,0] = 'Tanzania (United Republic)'
$geoData[123,0] = 'Thailand'
$geoData[124,0] = 'Timor-Leste'
$geoData[125,0] = 'Togo'
$geoData[126,0] = 'Tunisia'
$geoData[127,0] = 'Turkmenistan'
$geoData[128,0] = 'Uganda'
$geoData[129,0] = 'Ukraine'
$geoData[130,0] = 'Uzbekistan'
$geoData[131,0] = 'Venezuela (Bolivarian Republic)'
$geoData[132,0] = 'Viet Nam'
$geoData[133,0] = 'Zambia'
$geoData[134,0] = 'Zanzibar'
$geoData[135,0] = 'Zimbabwe'

$ws2.Range("A1:C136").Value = $geoData

# Leave the Geography sheet active, with B11 selected (matches the author's
# last on-screen selection when the join table was added).
$ws2.Activate()
$ws2.Range("B11").Select()
